# bookshelf.index.pptx — "add save and reload queue function" edit
#
# EMU -> point helper. PowerPoint COM's Shape.Left/Top/Width/Height are in
# points (1 pt = 914400/72 EMU); the interop layer truncates the EMU value
# it computes from the point we hand it, so nudge by +0.5 EMU before
# converting down to points to land exactly on the target integer EMU.
function EmuToPt([double]$emu) {
    return ($emu + 0.5) / 914400.0 * 72.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 ("./src/database/" flow)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Shape 6 "직사각형 5" — nudge position
$sh6 = $s4.Shapes.Item("직사각형 5")
$sh6.Left = EmuToPt 150470
$sh6.Top  = EmuToPt 191513

# Shape 7 "사각형: 둥근 모서리 6" (./init.js) — nudge position
$sh7 = $s4.Shapes.Item("사각형: 둥근 모서리 6")
$sh7.Left = EmuToPt 1597307
$sh7.Top  = EmuToPt 748317

# Shape 8 "직사각형 7" — reposition/resize and update text
#   "tables" -> "Tables," plus a new "Export database pointer" paragraph
$sh8 = $s4.Shapes.Item("직사각형 7")
$sh8.Left   = EmuToPt 4317357
$sh8.Top    = EmuToPt 490963
$sh8.Width  = EmuToPt 2708475
$sh8.Height = EmuToPt 923330
$sh8.TextFrame.TextRange.Text = "Create fundamental" + [char]13 + "Tables," + [char]13 + "Export database pointer"

# Connector 10 "직선 화살표 연결선 9" (sh7 -> sh8) — follows the moved shapes
$cxn10 = $s4.Shapes.Item("직선 화살표 연결선 9")
$cxn10.Left   = EmuToPt 3414532
$cxn10.Top    = EmuToPt 952628
$cxn10.Width  = EmuToPt 902825
$cxn10.Height = EmuToPt 1

# Shape 11 "사각형: 둥근 모서리 10" ("나머지 쿼리는 필요할 때마다 추가") — removed
$sh11 = $s4.Shapes.Item("사각형: 둥근 모서리 10")
$sh11.Delete()

# ---------------------------------------------------------------------
# Slide 5 (scan/meta -> scanner/agent flows)
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# Shape 11 "사각형: 둥근 모서리 10" (./scan.js) — widen box, fix filename
$s5sh11 = $s5.Shapes.Item("사각형: 둥근 모서리 10")
$s5sh11.Width  = EmuToPt 1470272
$s5sh11.Height = EmuToPt 408623
$s5sh11.TextFrame.TextRange.Text = "./scanner.js"

# Shape 12 "사각형: 둥근 모서리 11" (./meta.js) — widen box, rename (typo kept
# intentionally, matching the source commit: "agnet.js")
$s5sh12 = $s5.Shapes.Item("사각형: 둥근 모서리 11")
$s5sh12.Width  = EmuToPt 1470272
$s5sh12.Height = EmuToPt 408623
$s5sh12.TextFrame.TextRange.Text = "./agnet.js"

# Shape 105 "사각형: 둥근 모서리 104" (./scan.js) — widen box, fix filename
$s5sh105 = $s5.Shapes.Item("사각형: 둥근 모서리 104")
$s5sh105.Width  = EmuToPt 1510736
$s5sh105.Height = EmuToPt 408623
$s5sh105.TextFrame.TextRange.Text = "./scanner.js"

# Shape 106 "사각형: 둥근 모서리 105" (./meta.js) — widen box, rename
$s5sh106 = $s5.Shapes.Item("사각형: 둥근 모서리 105")
$s5sh106.Width  = EmuToPt 1510736
$s5sh106.Height = EmuToPt 408623
$s5sh106.TextFrame.TextRange.Text = "./agent.js"

# Connector 152 "직선 화살표 연결선 151" (105 -> 129) — shifts right with the
# widened box it starts from
$cxn152 = $s5.Shapes.Item("직선 화살표 연결선 151")
$cxn152.Left   = EmuToPt 8373431
$cxn152.Top    = EmuToPt 2441089
$cxn152.Width  = EmuToPt 855959
$cxn152.Height = EmuToPt 3735

# Connector 154 "직선 화살표 연결선 153" (106 -> 130) — shifts right with the
# widened box it starts from
$cxn154 = $s5.Shapes.Item("직선 화살표 연결선 153")
$cxn154.Left   = EmuToPt 8369529
$cxn154.Top    = EmuToPt 2978875
$cxn154.Width  = EmuToPt 859861
$cxn154.Height = EmuToPt 5346
